$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("RunManager")
$ws2 = $wb.Worksheets.Item("TestData")

# RunManager: rename row2's test to masterModuleTest, drop the loginTest/logoutTest rows
$ws1.Range("A2").Value = "masterModuleTest"
$ws1.Rows("3:4").Delete()

# Update selections on each sheet, selecting TestData first so RunManager
# ends up as the active/selected tab (matches the target file).
[void]$ws2.Range("F1:I1048576").Select()
[void]$ws1.Range("E1:H1048576").Select()
